# Trim stray leading/trailing whitespace from the text values in the
# "join*/majority/dissent/..." vote-description columns (B:J), fixing
# wiki-scrape typos as described in the commit message ("Fixing more
# wiki typos").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $trimmed = $val.Trim()
            if ($trimmed -ne $val) {
                $cell.Value = $trimmed
            }
        }
    }
}
